# Graph from results to gap %
#
# Updates the "Execution Time during development for each instances" bar
# chart (slide 20 / chart3.xml) so each series plots the "gap %" column
# instead of the old "results" column. The chart's data is backed by an
# external (not embedded) workbook, so the only writable surface exposed
# by the COM object model is Series.Values (the cached <c:numCache> point
# values) - the source-range formula text and its linked number format are
# not editable from this environment, so only the values are refreshed
# here.

$p = $ppt.ActivePresentation

# Locate the bar chart with 5 series (v0.1 / v0.3 / v0.4 / v0.5 / v0.7) -
# this is the "Execution Time during development for each instances" chart,
# normally on slide 20, but we search defensively instead of hardcoding the
# slide index.
$chart = $null
for ($si = 1; $si -le $p.Slides.Count -and $chart -eq $null; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $candidate = $s.Shapes.Item($i)
        if ($candidate.HasChart) {
            $candidateChart = $candidate.Chart
            if ($candidateChart.SeriesCollection().Count -eq 5) {
                $chart = $candidateChart
                break
            }
        }
    }
}

$sc = $chart.SeriesCollection()

$newValues = @{
    1 = @(12.07854352251514, 44.0556031579293, 68.57376207014956, 68.45483635130971, 40.59262994799748, 64.22637225914396, 49.24926129078125)
    2 = @(3.594325139842013, 29.76816843672586, 51.90334949015786, 49.03069791357022, 41.16622431430861, 51.87654906525632, 27.90661692619954)
    3 = @(1.047085514330305, 21.1248322472501, 30.42705220242081, 30.26742075690127, 26.48530904683112, 32.85913515764747, 21.99634617908459)
    4 = @(0.44084248345852, 15.79477493040673, 20.83367265188321, 19.97871767513666, 19.6796894033014, 23.03839809655065, 25.35942953350353)
    5 = @(0.0307369037511316, 7.932220110046956, 10.09705649676077, 8.912530733289074, 12.62602895271822, 19.26119153459033, 4.225734134728905)
}

for ($i = 1; $i -le $sc.Count; $i++) {
    $ser = $sc.Item($i)
    $ser.Values = $newValues[$i]
}
